$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Förändrad" (changed) date moved from 2023-09-20 (45189) to 2023-09-21 (45190)
# for every existing data row (2-157).
$ws.Range("C2:C157").Value = 45190

# Row 157 picks up an explicit row-height now that a new row follows it.
$ws.Rows.Item(157).RowHeight = 15

# New cleaning notice appended as row 158.
$ws.Cells.Item(158, 1).Value = "A 44224-2023"

$ws.Cells.Item(158, 2).Value = 45188
$ws.Cells.Item(158, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(158, 3).Value = 45190
$ws.Cells.Item(158, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(158, 4).Value = "ÖREBRO LÄN"
$ws.Cells.Item(158, 5).Value = "KARLSKOGA"
$ws.Cells.Item(158, 6).Value = "Sveaskog"
$ws.Cells.Item(158, 7).Value = 0.8
$ws.Cells.Item(158, 8).Value = 0
$ws.Cells.Item(158, 9).Value = 0
$ws.Cells.Item(158, 10).Value = 0
$ws.Cells.Item(158, 11).Value = 0
$ws.Cells.Item(158, 12).Value = 0
$ws.Cells.Item(158, 13).Value = 0
$ws.Cells.Item(158, 14).Value = 0
$ws.Cells.Item(158, 15).Value = 0
$ws.Cells.Item(158, 16).Value = 0
$ws.Cells.Item(158, 17).Value = 0

$ws.Cells.Item(158, 18).Value = ""
$ws.Cells.Item(158, 18).WrapText = $true
